$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 30 (shifts existing rows 30-161 down to 31-162)
$ws.Rows.Item(30).EntireRow.Insert()

# Populate the newly inserted row 30 with data
$ws.Range("A30").Value = 10
$ws.Range("B30").Value = "Vega Modelo de Temuco"
$ws.Range("C30").Value = "La Araucanía"
$ws.Range("D30").Value = 44764
$ws.Range("E30").Value = 9
$ws.Range("F30").Value = 100112012
$ws.Range("G30").Value = "Espinaca"
$ws.Range("H30").Value = "Sin especificar"
$ws.Range("I30").Value = "Primera"
$ws.Range("J30").Value = 85
$ws.Range("K30").Value = 15000
$ws.Range("L30").Value = 15000
$ws.Range("M30").Value = 15000
$ws.Range("N30").Value = "$/docena de atados"
$ws.Range("O30").Value = "Región de La Araucanía"
$ws.Range("P30").Value = 5000
$ws.Range("Q30").Value = 3
$ws.Range("R30").Value = "Hortaliza"
